# Adds "Start Date" (col K) / "End Date" (col L) values to the session
# schedule on the Template sheet. Sessions with an early start (8/25/2019 -
# 10/18/2019) get one pair of dates, sessions that start later
# (10/21/2019 - 12/13/2019) get the other pair - this is how the sheet
# flags a "2nd session" when the two date ranges diverge.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (Start Date serial, End Date serial)
$dates = [ordered]@{
    4  = @(43702, 43756)
    5  = @(43759, 43812)
    6  = @(43702, 43756)
    7  = @(43759, 43812)
    8  = @(43702, 43756)
    9  = @(43759, 43812)
    11 = @(43702, 43756)
    12 = @(43759, 43812)
    16 = @(43702, 43756)
    17 = @(43759, 43812)
    20 = @(43702, 43756)
    21 = @(43759, 43812)
    26 = @(43702, 43756)
    27 = @(43759, 43812)
    35 = @(43702, 43756)
    36 = @(43759, 43812)
    37 = @(43702, 43756)
    38 = @(43759, 43812)
}

# K5:L5 already carries the date number format used everywhere else in the
# sheet - copy its formatting onto every K/L cell we're about to fill so the
# new cells match the existing date columns instead of picking up a brand
# new style.
$formatSource = $ws.Range("K5:L5")

foreach ($r in $dates.Keys) {
    $pair = $dates[$r]
    $target = $ws.Range("K" + $r + ":L" + $r)

    $formatSource.Copy()
    $target.PasteSpecial(-4122)

    $ws.Cells.Item($r, 11).Value = $pair[0]
    $ws.Cells.Item($r, 12).Value = $pair[1]
}

$excel.CutCopyMode = 0

# Restore the scroll position / selection left behind by the edit.
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("K26:L27").Select()
